$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# A new (blank) column is inserted before column N, shifting the old
# N/O/P columns (Late / Heading / Outstanding) one column to the right.
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab (it was "Transactions"
# before), restore the horizontal scroll position and leave the cursor
# on the far right of the sheet, as was left by the author.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("U11").Select()
